$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Cells.Item(2, 1).Value2 = 'Payton Pritchard'
$ws.Cells.Item(2, 2).Value2 = 'PG'
$ws.Cells.Item(2, 3).Value2 = 'Boston Celtics'
$ws.Cells.Item(3, 1).Value2 = 'Donovan Mitchell'
$ws.Cells.Item(3, 2).Value2 = 'PG,SG'
$ws.Cells.Item(3, 3).Value2 = 'Cleveland Cavaliers'
$ws.Cells.Item(4, 1).Value2 = 'Jamal Murray'
$ws.Cells.Item(4, 2).Value2 = 'PG,SG'
$ws.Cells.Item(4, 3).Value2 = 'Denver Nuggets'
$ws.Cells.Item(5, 1).Value2 = 'Marcus Smart'
$ws.Cells.Item(5, 2).Value2 = 'PG,SG'
$ws.Cells.Item(5, 3).Value2 = 'Memphis Grizzlies'
$ws.Cells.Item(6, 1).Value2 = 'Josh Hart'
$ws.Cells.Item(6, 2).Value2 = 'SF,PF'
$ws.Cells.Item(6, 3).Value2 = 'New York Knicks'
$ws.Cells.Item(7, 1).Value2 = 'Tari Eason'
$ws.Cells.Item(7, 2).Value2 = 'SF,PF'
$ws.Cells.Item(7, 3).Value2 = 'Houston Rockets'
$ws.Cells.Item(8, 1).Value2 = 'Michael Porter Jr.'
$ws.Cells.Item(8, 2).Value2 = 'SF,PF'
$ws.Cells.Item(8, 3).Value2 = 'Denver Nuggets'
$ws.Cells.Item(9, 1).Value2 = 'Kristaps Porzingis'
$ws.Cells.Item(9, 2).Value2 = 'PF,C'
$ws.Cells.Item(9, 3).Value2 = 'Boston Celtics'
$ws.Cells.Item(10, 1).Value2 = 'Myles Turner'
$ws.Cells.Item(10, 2).Value2 = 'C'
$ws.Cells.Item(10, 3).Value2 = 'Indiana Pacers'
$ws.Cells.Item(11, 1).Value2 = 'Domantas Sabonis'
$ws.Cells.Item(11, 2).Value2 = 'C'
$ws.Cells.Item(11, 3).Value2 = 'Sacramento Kings'
$ws.Cells.Item(12, 1).Value2 = 'Deandre Ayton'
$ws.Cells.Item(12, 2).Value2 = 'C'
$ws.Cells.Item(12, 3).Value2 = 'Portland Trail Blazers'
$ws.Cells.Item(13, 1).Value2 = 'Victor Wembanyama'
$ws.Cells.Item(13, 2).Value2 = 'C'
$ws.Cells.Item(13, 3).Value2 = 'San Antonio Spurs'
$ws.Cells.Item(14, 1).Value2 = 'Dyson Daniels'
$ws.Cells.Item(14, 2).Value2 = 'PG,SG'
$ws.Cells.Item(14, 3).Value2 = 'Atlanta Hawks'
$ws.Cells.Item(15, 1).Value2 = 'Malik Beasley'
$ws.Cells.Item(15, 2).Value2 = 'SG'
$ws.Cells.Item(15, 3).Value2 = 'Detroit Pistons'
$ws.Cells.Item(16, 1).Value2 = 'De''Andre Hunter'
$ws.Cells.Item(16, 2).Value2 = 'SF,PF'
$ws.Cells.Item(16, 3).Value2 = 'Atlanta Hawks'
$ws.Cells.Item(17, 1).Value2 = 'Zion Williamson'
$ws.Cells.Item(17, 2).Value2 = 'PF,C'
$ws.Cells.Item(17, 3).Value2 = 'New Orleans Pelicans'
$ws.Cells.Item(18, 1).Value2 = 'Bradley Beal'
$ws.Cells.Item(18, 2).Value2 = 'PG,SG,SF'
$ws.Cells.Item(18, 3).Value2 = 'Phoenix Suns'
